# Allocation rule updated with 5 and 10 mi radius columns
$wb = $excel.ActiveWorkbook

$meansWs = $wb.Worksheets.Item("Means")
$sdWs = $wb.Worksheets.Item("Standard Deviations")

# --- Header row (row 1): add columns F and G ---
$meansWs.Cells.Item(1, 6).Value = "Within 5 miles of HFC production facility"
$meansWs.Cells.Item(1, 7).Value = "Within 10 miles of HFC production facility"

$sdWs.Cells.Item(1, 6).Value = "Within 5 mile of HFC production facility SD"
$sdWs.Cells.Item(1, 7).Value = "Within 10 mile of HFC production facility SD"

# --- Means sheet data updates ---
# Row 2: % White
$meansWs.Cells.Item(2, 6).Value = 82
$meansWs.Cells.Item(2, 7).Value = 62

# Row 3: % Black or African American
$meansWs.Cells.Item(3, 6).Value = 15
$meansWs.Cells.Item(3, 7).Value = 35

# Row 4: % Other
$meansWs.Cells.Item(4, 6).Value = 2.9
$meansWs.Cells.Item(4, 7).Value = 3.4

# Row 5: % Hispanic
$meansWs.Cells.Item(5, 6).Value = 3.4
$meansWs.Cells.Item(5, 7).Value = 4.5

# Row 6: Median Income [1,000 2019$]
$meansWs.Cells.Item(6, 6).Value = 54
$meansWs.Cells.Item(6, 7).Value = 45

# Row 7: % Below Poverty Line
$meansWs.Cells.Item(7, 6).Value = 11
$meansWs.Cells.Item(7, 7).Value = 13

# Row 8: % Below Half the Poverty Line
$meansWs.Cells.Item(8, 6).Value = 4.2
$meansWs.Cells.Item(8, 7).Value = 7.7

# Row 9: Total Cancer Risk (per million) -- existing B-E values also change
$meansWs.Cells.Item(9, 2).Value = 26
$meansWs.Cells.Item(9, 3).Value = 34
$meansWs.Cells.Item(9, 4).Value = 50
$meansWs.Cells.Item(9, 5).Value = 50
$meansWs.Cells.Item(9, 6).Value = 50
$meansWs.Cells.Item(9, 7).Value = 49

# Row 10: Total Respiratory (hazard quotient) -- existing B-E values also change
$meansWs.Cells.Item(10, 2).Value = 0.32
$meansWs.Cells.Item(10, 3).Value = 0.47
$meansWs.Cells.Item(10, 4).Value = 0.6
$meansWs.Cells.Item(10, 5).Value = 0.6
$meansWs.Cells.Item(10, 6).Value = 0.57
$meansWs.Cells.Item(10, 7).Value = 0.54

# --- Standard Deviations sheet data updates ---
# Row 2: % White
$sdWs.Cells.Item(2, 6).Value = 15
$sdWs.Cells.Item(2, 7).Value = 31

# Row 3: % Black or African American
$sdWs.Cells.Item(3, 6).Value = 17
$sdWs.Cells.Item(3, 7).Value = 33

# Row 4: % Other
$sdWs.Cells.Item(4, 6).Value = 4.8
$sdWs.Cells.Item(4, 7).Value = 4.4

# Row 5: % Hispanic
$sdWs.Cells.Item(5, 6).Value = 2.2
$sdWs.Cells.Item(5, 7).Value = 7.6

# Row 6: Median Income [1,000 2019$]
$sdWs.Cells.Item(6, 6).Value = 30
$sdWs.Cells.Item(6, 7).Value = 21

# Row 7: % Below Poverty Line
$sdWs.Cells.Item(7, 6).Value = 11
$sdWs.Cells.Item(7, 7).Value = 14

# Row 8: % Below Half the Poverty Line
$sdWs.Cells.Item(8, 6).Value = 2.3
$sdWs.Cells.Item(8, 7).Value = 9.7

# Row 9: Total Cancer Risk (per million) -- existing B-E values also change
$sdWs.Cells.Item(9, 2).Value = 8.6
$sdWs.Cells.Item(9, 3).Value = 5.5
$sdWs.Cells.Item(9, 4).Value = 0
$sdWs.Cells.Item(9, 5).Value = 0
$sdWs.Cells.Item(9, 6).Value = 0
$sdWs.Cells.Item(9, 7).Value = 4.9

# Row 10: Total Respiratory (hazard quotient) -- existing B-E values also change
$sdWs.Cells.Item(10, 2).Value = 0.14
$sdWs.Cells.Item(10, 3).Value = 0.056
$sdWs.Cells.Item(10, 4).Value = 0
$sdWs.Cells.Item(10, 5).Value = 0
$sdWs.Cells.Item(10, 6).Value = 0.052
$sdWs.Cells.Item(10, 7).Value = 0.064
